{"js": "// Update the date line and the 25 division answers in the practice table.\n// Every text run keeps its original run/paragraph formatting; only the\n// literal text content changes, matching the source diff.\n\nconst titleParas = context.document.body.paragraphs;\ntitleParas.load(\"items\");\nawait context.sync();\n// The very first paragraph in the body holds the date heading.\ntitleParas.items[0].insertText(\"2024-03-02 Saturday\", \"Replace\");\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\nconst table = tables.items[0];\n\n// Map of table row index -> new cell values (left to right).\n// Only rows 0, 4, 8, 12, 16 contain answers; the other rows are spacer rows.\nconst rowUpdates = {\n  0: [\"94\u00f79=10, 4\", \"62\u00f73=20, 2\", \"83\u00f73=27, 2\", \"70\u00f74=17, 2\", \"98\u00f75=19, 3\"],\n  4: [\"23\u00f75=4, 3\", \"31\u00f73=10, 1\", \"65\u00f73=21, 2\", \"16\u00f72=8, 0\", \"16\u00f72=8, 0\"],\n  8: [\"65\u00f77=9, 2\", \"44\u00f76=7, 2\", \"48\u00f77=6, 6\", \"46\u00f74=11, 2\", \"33\u00f79=3, 6\"],\n  12: [\"98\u00f79=10, 8\", \"72\u00f76=12, 0\", \"11\u00f78=1, 3\", \"31\u00f73=10, 1\", \"13\u00f78=1, 5\"],\n  16: [\"21\u00f72=10, 1\", \"53\u00f74=13, 1\", \"47\u00f72=23, 1\", \"72\u00f75=14, 2\", \"24\u00f79=2, 6\"],\n};\n\nfor (const rowIndexStr of Object.keys(rowUpdates)) {\n  const rowIndex = parseInt(rowIndexStr, 10);\n  const values = rowUpdates[rowIndexStr];\n  for (let col = 0; col < values.length; col++) {\n    const cell = table.getCell(rowIndex, col);\n    cell.body.paragraphs.load(\"items\");\n    await context.sync();\n    cell.body.paragraphs.items[0].insertText(values[col], \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the date line and the 25 division answers in the practice table.\n# Assigning to Range.Text only swaps the literal characters, so the existing\n# run/paragraph formatting (fonts, size, alignment) is left untouched.\n\n$d = $word.ActiveDocument\n\n# The first paragraph in the document holds the date heading.\n$d.Paragraphs.Item(1).Range.Text = \"2024-03-02 Saturday\"\n\n$table = $d.Tables.Item(1)\n\n# Row (1-based) -> new cell values left to right. Only these five rows hold\n# answers; the rows in between are blank spacer rows and are left alone.\n$rowUpdates = @{\n    1  = @(\"94\u00f79=10, 4\", \"62\u00f73=20, 2\", \"83\u00f73=27, 2\", \"70\u00f74=17, 2\", \"98\u00f75=19, 3\")\n    5  = @(\"23\u00f75=4, 3\", \"31\u00f73=10, 1\", \"65\u00f73=21, 2\", \"16\u00f72=8, 0\", \"16\u00f72=8, 0\")\n    9  = @(\"65\u00f77=9, 2\", \"44\u00f76=7, 2\", \"48\u00f77=6, 6\", \"46\u00f74=11, 2\", \"33\u00f79=3, 6\")\n    13 = @(\"98\u00f79=10, 8\", \"72\u00f76=12, 0\", \"11\u00f78=1, 3\", \"31\u00f73=10, 1\", \"13\u00f78=1, 5\")\n    17 = @(\"21\u00f72=10, 1\", \"53\u00f74=13, 1\", \"47\u00f72=23, 1\", \"72\u00f75=14, 2\", \"24\u00f79=2, 6\")\n}\n\nforeach ($rowIndex in $rowUpdates.Keys) {\n    $values = $rowUpdates[$rowIndex]\n    for ($col = 1; $col -le $values.Length; $col++) {\n        $table.Cell($rowIndex, $col).Range.Text = $values[$col - 1]\n    }\n}\n"}
